# Update vm_pu.xlsx res_bus values for the 380 kV case (Case_0_6)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.007233763412072
$ws.Range("D2").Value = 1.009913254666529
$ws.Range("E2").Value = 1.009822179915484
$ws.Range("F2").Value = 1.005418171149548
$ws.Range("I2").Value = 1.023594999628091
$ws.Range("J2").Value = 1.012506593353203
$ws.Range("K2").Value = 1.012785033388607
$ws.Range("L2").Value = 1.012694234278072
$ws.Range("M2").Value = 1.008303619376536
$ws.Range("N2").Value = 1.00796475591635
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.009287095432226
$ws.Range("D3").Value = 1.01176408710081
$ws.Range("E3").Value = 1.011604225035834
$ws.Range("F3").Value = 1.008140331328463
$ws.Range("I3").Value = 1.023504579208683
$ws.Range("J3").Value = 1.014185389888849
$ws.Range("K3").Value = 1.014436890927128
$ws.Range("L3").Value = 1.01427747657216
$ws.Range("M3").Value = 1.010823321684793
$ws.Range("N3").Value = 1.008557101822973
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.010609914077341
$ws.Range("D4").Value = 1.012956425105847
$ws.Range("E4").Value = 1.012752198477128
$ws.Range("F4").Value = 1.009894655252362
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.015265863463752
$ws.Range("K4").Value = 1.015500096321514
$ws.Range("L4").Value = 1.015296411973174
$ws.Range("M4").Value = 1.012446481218014
$ws.Range("N4").Value = 1.008937390812613
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.011164667208386
$ws.Range("D5").Value = 1.013456451286344
$ws.Range("E5").Value = 1.013233607139253
$ws.Range("F5").Value = 1.010630526950308
$ws.Range("I5").Value = 1.023416932628352
$ws.Range("J5").Value = 1.015718730223423
$ws.Range("K5").Value = 1.015945741712944
$ws.Range("L5").Value = 1.015723475712681
$ws.Range("M5").Value = 1.013127169341356
$ws.Range("N5").Value = 1.009096556300793
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.011257733956607
$ws.Range("D6").Value = 1.013540336460009
$ws.Range("E6").Value = 1.013314368262514
$ws.Range("F6").Value = 1.01075398803822
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.015794689376722
$ws.Range("K6").Value = 1.016020490579834
$ws.Range("L6").Value = 1.015795106326709
$ws.Range("M6").Value = 1.013241362293345
$ws.Range("N6").Value = 1.009123239652299
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.01061733201275
$ws.Range("D7").Value = 1.012963111283719
$ws.Range("E7").Value = 1.012758635757744
$ws.Range("F7").Value = 1.009904494403763
$ws.Range("I7").Value = 1.02344298551034
$ws.Range("J7").Value = 1.015271920013607
$ws.Range("K7").Value = 1.015506056230265
$ws.Range("L7").Value = 1.015302123479047
$ws.Range("M7").Value = 1.01245558317387
$ws.Range("N7").Value = 1.008939520358412
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.007928923536864
$ws.Range("D8").Value = 1.010539862416142
$ws.Range("E8").Value = 1.010425509666386
$ws.Range("F8").Value = 1.006339639084466
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.01307517316777
$ws.Range("K8").Value = 1.013344474725078
$ws.Range("L8").Value = 1.013230459314929
$ws.Range("M8").Value = 1.009156696789249
$ws.Range("N8").Value = 1.008165568793856
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.003145426914804
$ws.Range("D9").Value = 1.006228032427491
$ws.Range("E9").Value = 1.006273668063492
$ws.Range("F9").Value = 1.000001172251527
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.009158333178679
$ws.Range("K9").Value = 1.009490890245721
$ws.Range("L9").Value = 1.009536367171993
$ws.Range("M9").Value = 1.003285835256702
$ws.Range("N9").Value = 1.006778356089766
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 0.9999231654433728
$ws.Range("D10").Value = 1.003323481638483
$ws.Range("E10").Value = 1.003476664060879
$ws.Range("F10").Value = 0.9957340015561771
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.00651439031437
$ws.Range("K10").Value = 1.006890041508089
$ws.Range("L10").Value = 1.007042636510944
$ws.Range("M10").Value = 0.9993299166240076
$ws.Range("N10").Value = 1.005837174268722
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 0.9985195009428353
$ws.Range("D11").Value = 1.002058230093306
$ws.Range("E11").Value = 1.002258213239049
$ws.Range("F11").Value = 0.9938756573484834
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.005361353862241
$ws.Range("K11").Value = 1.005755897172605
$ws.Range("L11").Value = 1.005955082248305
$ws.Range("M11").Value = 0.9976062853082813
$ws.Range("N11").Value = 1.005425605319509
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 0.9979968113523644
$ws.Range("D12").Value = 1.001587085480179
$ws.Range("E12").Value = 1.001804489229526
$ws.Range("F12").Value = 0.9931837220361893
$ws.Range("I12").Value = 1.023938164268129
$ws.Range("J12").Value = 1.004931797492527
$ws.Range("K12").Value = 1.005333394152455
$ws.Range("L12").Value = 1.005549917675707
$ws.Range("M12").Value = 0.9969643836886319
$ws.Range("N12").Value = 1.005272112154014
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 0.998108989854532
$ws.Range("D13").Value = 1.001688201332467
$ws.Range("E13").Value = 1.00190186661117
$ws.Range("F13").Value = 0.9933322208667887
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.005023996734401
$ws.Range("K13").Value = 1.005424078780853
$ws.Range("L13").Value = 1.005636881635679
$ws.Range("M13").Value = 0.9971021502599321
$ws.Range("N13").Value = 1.00530506513454
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 0.9984763221568785
$ws.Range("D14").Value = 1.002019309300045
$ws.Range("E14").Value = 1.002220731678253
$ws.Range("F14").Value = 0.9938184960838859
$ws.Range("I14").Value = 1.023922712353276
$ws.Range("J14").Value = 1.005325872665321
$ws.Range("K14").Value = 1.005720998267353
$ws.Range("L14").Value = 1.005921615865151
$ws.Range("M14").Value = 0.997553259913063
$ws.Range("N14").Value = 1.005412930209644
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 0.9987024733256592
$ws.Range("D15").Value = 1.002223159157142
$ws.Range("E15").Value = 1.00241704320311
$ws.Range("F15").Value = 0.9941178838525467
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.005511699354268
$ws.Range("K15").Value = 1.005903775928103
$ws.Range("L15").Value = 1.006096890186832
$ws.Range("M15").Value = 0.9978309805352872
$ws.Range("N15").Value = 1.005479307154862
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.000016141277715
$ws.Range("D16").Value = 1.003407289689546
$ws.Range("E16").Value = 1.003557371151442
$ws.Range("F16").Value = 0.9958571041249777
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.006590737850094
$ws.Range("K16").Value = 1.006965140220006
$ws.Range("L16").Value = 1.007114647702511
$ws.Range("M16").Value = 0.9994440778352311
$ws.Range("N16").Value = 1.005864402752942
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.000837890047495
$ws.Range("D17").Value = 1.004148012526777
$ws.Range("E17").Value = 1.0042706812522
$ws.Range("F17").Value = 0.9969451783998167
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.007265371433142
$ws.Range("K17").Value = 1.007628750151268
$ws.Range("L17").Value = 1.007750960917956
$ws.Range("M17").Value = 1.000453025260634
$ws.Range("N17").Value = 1.006104875074253
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.001316394315354
$ws.Range("D18").Value = 1.004579336601008
$ws.Range("E18").Value = 1.004686038224454
$ws.Range("F18").Value = 0.9975788125116317
$ws.Range("I18").Value = 1.023826118601784
$ws.Range("J18").Value = 1.007658085910057
$ws.Range("K18").Value = 1.008015057034474
$ws.Range("L18").Value = 1.008121365968804
$ws.Range("M18").Value = 1.001040499793651
$ws.Range("N18").Value = 1.006244750435481
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.001479416177004
$ws.Range("D19").Value = 1.004726284779849
$ws.Range("E19").Value = 1.004827545758709
$ws.Range("F19").Value = 0.99779469414251
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.007791858856404
$ws.Range("K19").Value = 1.00814664890886
$ws.Range("L19").Value = 1.008247539050799
$ws.Range("M19").Value = 1.001240641009048
$ws.Range("N19").Value = 1.006292378869608
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.000749808057148
$ws.Range("D20").Value = 1.004068615447174
$ws.Range("E20").Value = 1.004194222981725
$ws.Range("F20").Value = 0.9968285444312279
$ws.Range("I20").Value = 1.023846096030143
$ws.Range("J20").Value = 1.007193071390838
$ws.Range("K20").Value = 1.007557630526282
$ws.Range("L20").Value = 1.007682767903857
$ws.Range("M20").Value = 1.000344881465634
$ws.Range("N20").Value = 1.0060791149003
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 0.9983681883574533
$ws.Range("D21").Value = 1.001921838967304
$ws.Range("E21").Value = 1.002126865545669
$ws.Range("F21").Value = 0.9936753466134768
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.005237013008493
$ws.Range("K21").Value = 1.005633597157727
$ws.Range("L21").Value = 1.005837802090324
$ws.Range("M21").Value = 0.9974204659168552
$ws.Range("N21").Value = 1.005381183798786
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 0.9968631906893453
$ws.Range("D22").Value = 1.000565264582737
$ws.Range("E22").Value = 1.000820438517567
$ws.Range("F22").Value = 0.9916831474676434
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.003999810658025
$ws.Range("K22").Value = 1.004416738918358
$ws.Range("L22").Value = 1.004670847506421
$ws.Range("M22").Value = 0.9955720865436296
$ws.Range("N22").Value = 1.004938784099908
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 0.9976617515572058
$ws.Range("D23").Value = 1.001285068720396
$ws.Range("E23").Value = 1.001513637737121
$ws.Range("F23").Value = 0.9927401879471101
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.00465638420563
$ws.Range("K23").Value = 1.005062507507653
$ws.Range("L23").Value = 1.005290142572135
$ws.Range("M23").Value = 0.9965528862430495
$ws.Range("N23").Value = 1.005173652397911
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.000789611011589
$ws.Range("D24").Value = 1.004104493814723
$ws.Range("E24").Value = 1.004228773359168
$ws.Range("F24").Value = 0.9968812494534719
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.007225743100769
$ws.Range("K24").Value = 1.007589768787545
$ws.Range("L24").Value = 1.007713583694289
$ws.Range("M24").Value = 1.000393750163363
$ws.Range("N24").Value = 1.006090756012724
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.00438777342681
$ws.Range("D25").Value = 1.007347887935682
$ws.Range("E25").Value = 1.007352013251211
$ws.Range("F25").Value = 1.001646885398627
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.010176558077594
$ws.Range("K25").Value = 1.010492603475836
$ws.Range("L25").Value = 1.010496715010997
$ws.Range("M25").Value = 1.004810763952433
$ws.Range("N25").Value = 1.007139822571722
